$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    27  = "0fa840aedb820ecb2432ca2597b63195"
    32  = "38807853720eb0201246d26fcfc87872"
    194 = "016fd01c9c0b9c0613409c3bd1823181"
    225 = "7fb50065b239688f7187b753c4a0fecf"
    286 = "c7adf4416ce8b934e1da9eb2b284a86e"
    298 = "62b2c583fda5fcaae2c96116de99e53a"
    303 = "0b8528a062afeffbab4736b6e43d8ee4"
    356 = "2b9bd6197117b0227f88ee1389b6fd3c"
    364 = "884ca1dc8b42bf2cf75c1180c6035b0c"
    373 = "90954109c5d442f2adf8575dd44df35d"
    403 = "9bcbb2f255b19c6d0539e8cbda99fd05"
    409 = "235eaa3861f9fcefcf38b0240fa98e23"
    422 = "83915188c396798e95a3c5d842cea75a"
    426 = "930e9bd628ccd09c643cd2b4a4b8cfad"
    428 = "abd8e9546a63d64fa7b638163f0d318d"
    483 = "2e1cb344bc4c52acae679a96ac27c388"
    517 = "19ddbcf717fba0b769c5f1391bea6f6e"
    518 = "2effef6c964dcb44f3b44ffedd192277"
    533 = "379d4952f7b707ba2ab09e2a82834f06"
    553 = "869d246e47abf3ed3b32866c54f6b704"
    572 = "fd4b65761d7f17e0ff910e8e0c028f25"
    585 = "b990d52701b34f357d5d0ed1c21643ca"
    593 = "42ccffc015f83aa5688c9de71fc231df"
    597 = "ab05b9db032c806c05c33ff66bf55bda"
    747 = "2802ab1063279d54146223f696f20eb3"
    750 = "fb6579275369feca2249f6a62946d497"
    755 = "45cce2fdc22e2cfd7fa5302a2e549dab"
    798 = "a55cd841ae22817e9d8e75e6235c48a7"
    812 = "13a58843e6c5aa259a24eb42e4c67fbc"
    907 = "402455665a127b9c85bbd42f906ecf74"
    915 = "ab434232503911719a501da4bb02d3c9"
    945 = "0739e4252751d56b83824b70b671b54d"
    955 = "81016273d6b4ba3f8e0fca8df6f28510"
    973 = "f74152985e919acc18d24beda95e9ea2"
    977 = "2775c1aca94ef3be5e1cb93e632b9c32"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 2).Value = $updates[$row]
}
